$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking values
# (e.g. "582.03") are not auto-converted to numbers, matching the source
# workbook which stores these as inline strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '64.510.18'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '2.636.00'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '582.03'
$ws.Range('E5').Value = '  -2.47%  '
$ws.Range('D6').Value = '156.93'
$ws.Range('E6').Value = '  +0.68%  '
$ws.Range('D7').Value = '0.646'
$ws.Range('E7').Value = '  +3.09%  '
$ws.Range('E9').Value = '  -3.46%  '
$ws.Range('D10').Value = '5.82'
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('D11').Value = '0.391'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '28.70'
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').Value = '0.0000188'
$ws.Range('E14').Value = '  -4.11%  '
$ws.Range('D15').Value = '3.113.20'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').Value = '64.309.95'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').Value = '2.627.91'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '12.28'
$ws.Range('E18').Value = '  -2.97%  '
$ws.Range('D19').Value = '4.69'
$ws.Range('E19').Value = '  -1.61%  '
$ws.Range('D20').Value = '7.45'
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').Value = '347.13'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '68.28'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('D24').Value = '0.0000113'
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').Value = '1.76'
$ws.Range('E25').Value = '  +5.90%  '
$ws.Range('D26').Value = '9.46'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('D27').Value = '588.78'
$ws.Range('E27').Value = '  +9.49%  '
$ws.Range('E28').Value = '  +0.80%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '8.06'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '0.162'
$ws.Range('E30').Value = '  -1.93%  '
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('D33').Value = '6.71'
$ws.Range('E33').Value = '  +4.57%  '
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').Value = '5.37'
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('D36').Value = '0.413'
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('D37').Value = '20.01'
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').Value = '1.94'
$ws.Range('E39').Value = '  +1.36%  '
$ws.Range('D40').Value = '155.03'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').Value = '2.45'
$ws.Range('E42').Value = '  +6.15%  '
$ws.Range('D43').Value = '158.51'
$ws.Range('E43').Value = '  -1.75%  '
$ws.Range('D44').Value = '4.01'
$ws.Range('E44').Value = '  -1.50%  '
$ws.Range('D45').Value = '23.48'
$ws.Range('E45').Value = '  +4.28%  '
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('D47').Value = '0.637'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.102'
$ws.Range('E48').Value = '  +3.02%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = '0.0254'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').Value = '19.27'
$ws.Range('E50').Value = '  -2.15%  '
$ws.Range('D51').Value = '0.0₆0236'
$ws.Range('E51').Value = '  -5.01%  '
